$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows at the bottom (rows 40-43), matching the shared-string insertion order
$ws.Range("A40").Value = "Accupuncture"
$ws.Range("B40").Value = "health clinic"
$ws.Range("C40").Value = 2

$ws.Range("A41").Value = "Behavioral"
$ws.Range("B41").Value = "behavioral health clinic"
$ws.Range("C41").Value = 1

$ws.Range("A42").Value = "Government"
$ws.Range("B42").Value = "govt, not ph"
$ws.Range("C42").Value = 2

$ws.Range("A43").Value = "Outreach"
$ws.Range("B43").Value = "community based organization"
$ws.Range("C43").Value = 2

# Now re-sort the whole data range (A2:C43) ascending by column A to match final layout
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2:A43"))
$ws.Sort.SetRange($ws.Range("A2:C43"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Perform the secondary sort noted in the diff's remaining sortState (A3:C24) - this is the trailing
# recorded sort-state left in the file after the user's last explicit Sort action on that sub-range.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A3:A24"))
$ws.Sort.SetRange($ws.Range("A3:C24"))
$ws.Sort.Header = -4142
$ws.Sort.Apply()

# Update the active selection to match the final state
$ws.Range("C37").Select() | Out-Null

Write-Host "Added 4 new tier-default rows and re-sorted the lookup table."
